# Improve image data formats figures
#
# 1) Refresh the cached "datetimeFigureOut" field text (Date placeholder)
#    on the slide master and every slide layout: 16.09.24 -> 14.11.24
# 2) Slide 3: "Data type (bit-depth)" -> "Pixel data type (bit-depth)"
# 3) Slide 3: reposition/resize the "Straight Arrow Connector 18" connector

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text refresh (master + all custom layouts)
# ---------------------------------------------------------------------
$newDate = "14.11.24"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -ne $newDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 3 text edit: "Data type (bit-depth)" -> "Pixel data type (bit-depth)"
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$tb14 = $slide3.Shapes.Item("TextBox 14")
$tr14 = $tb14.TextFrame.TextRange
for ($i = 1; $i -le $tr14.Paragraphs().Count; $i++) {
    $para = $tr14.Paragraphs($i)
    $ptext = $para.Text.TrimEnd([char]13)
    if ($ptext -eq "Data type (bit-depth)") {
        $sub = $tr14.Characters($para.Start, $para.Length)
        $sub.Text = "Pixel data type (bit-depth)"
        break
    }
}

# ---------------------------------------------------------------------
# 3) Slide 3: move/resize "Straight Arrow Connector 18"
# ---------------------------------------------------------------------
$conn = $slide3.Shapes.Item("Straight Arrow Connector 18")
$conn.Left = 369.908081
$conn.Top = 126.495476
$conn.Width = 135.803894
$conn.Height = 44.814922
